$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell J1: trafo_id -> gridnode_id
$ws.Range("J1").Value = "gridnode_id"

# Update the active selection to match the recorded view state
$ws.Range("E6").Select()
